# Update "想去人数" (want-to-go count) values in column F for the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet.
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 305
$wsExpo.Range("F3").Value = 1224
$wsExpo.Range("F4").Value = 16931
$wsExpo.Range("F6").Value = 1656
$wsExpo.Range("F8").Value = 8
$wsExpo.Range("F9").Value = 5
$wsExpo.Range("F13").Value = 11735
$wsExpo.Range("F15").Value = 5
$wsExpo.Range("F16").Value = 1419
$wsExpo.Range("F17").Value = 4660
$wsExpo.Range("F18").Value = 472
$wsExpo.Range("F19").Value = 13
$wsExpo.Range("F21").Value = 71
$wsExpo.Range("F22").Value = 902

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 305
$wsAll.Range("F4").Value = 1224
$wsAll.Range("F5").Value = 16931
$wsAll.Range("F7").Value = 1656
$wsAll.Range("F9").Value = 8
$wsAll.Range("F10").Value = 5
$wsAll.Range("F16").Value = 11735
$wsAll.Range("F18").Value = 5
$wsAll.Range("F19").Value = 1419
$wsAll.Range("F20").Value = 4660
$wsAll.Range("F21").Value = 472
$wsAll.Range("F22").Value = 13
$wsAll.Range("F24").Value = 71
$wsAll.Range("F25").Value = 902
